$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Query text blocks ---
$qStat = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Belgian Malinois']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@

$qCases = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS ``Case ID``,
       coalesce(s.clinical_study_designation, '') AS ``Study Code``,
       coalesce(s.clinical_study_type, '') AS  ``Study Type``,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS ``Stage of Disease``,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS ``Weight (kg)``,
       coalesce(diag.best_response, '') AS ``Response to Treatment``,
       coalesce(co.cohort_description, '') AS ``Cohort``
Order by c.case_id LIMIT 100
"@

$qSamples = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Belgian Malinois'] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, 
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS ``Sample Site``,
        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,
        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,
        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,
        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,
        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,
        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,
        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``
Order by samp.sample_id LIMIT 100
"@

$qFiles = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois'] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS ``File Name``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_type, '') AS ``File Type``,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(samp.sample_id, '') AS ``Sample ID``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
Order By f.file_name LIMIT 100
"@

$qStudyFiles = @"
  MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE demo.breed  IN ['Belgian Malinois'] 
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS ``File Name``,
  coalesce(f.file_type, '') AS ``File Type``,
  coalesce("study", '') AS ``Association``,
  coalesce(f.file_description, '') AS ``Description``,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS ``Study Code``
Order By f.file_name LIMIT 100
"@

# --- Row 2 (CasesTab): update Cases query (B2) and StatQuery (C2) ---
$ws.Range("B2").Value = $qCases
$ws.Range("C2").Value = $qStat

# --- Row 3 (SamplesTab): fix B3 to correct Samples query, update C3 ---
$ws.Range("B3").Value = $qSamples
$ws.Range("C3").Value = $qStat

# --- Row 4 (FilesTab): fix B4 to the corrected Files query, update C4 ---
$ws.Range("B4").Value = $qFiles
$ws.Range("C4").Value = $qStat

# --- Row 5 (new StudyFilesTab) ---
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $qStudyFiles
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = $qStat
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = "TC06_Canine_Filter_Breed-BelgMalin_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC06_Canine_Filter_Breed-BelgMalin_WebData.xlsx"

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 302.39999999999998
$ws.Rows.Item(4).RowHeight = 220.2
$ws.Rows.Item(5).RowHeight = 282.60000000000002

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 91.6
$ws.Columns.Item(5).ColumnWidth = 44.3

# --- Sheet view ---
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("B1").Select()
